$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.855.93"
$ws.Range("D3").Value = "1.813.00"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'310.31"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.4640"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.3695"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").Value = "'0.07349"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "'0.8691"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "1.871.79"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").Value = "'5.337"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'6.507"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").Value = "'91.59"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'0.000008711"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "'1.002"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "'14.68"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "26.905.23"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "'5.342"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'10.54"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "2.124.12"
$ws.Range("E24").Value = "  +3.57%  "
$ws.Range("D25").Value = "'1.896"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").Value = "'151.91"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'18.36"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").Value = "'2.125"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("D29").Value = "'5.296"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "'0.08907"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'0.7568"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("D33").Value = "'1.150"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("D34").Value = "'2.929"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'4.457"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").Value = "'1.001"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'1.094"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "'2.947"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5338"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.209"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.364"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "'0.1658"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").Value = "'8.429"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").Value = "'0.4936"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").Value = "'10.31"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").Value = "'1.001"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'102.95"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "'0.06265"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
